# Shear Center Reference point added
# Update the data table (D2, E2) and move the active selection to V7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values: D2 2 -> 1, E2 3 -> 2
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2

# Move the active selection to V7 (matches the sheetView selection change)
$ws.Activate()
$ws.Range("V7").Select()
